$wb = $excel.ActiveWorkbook

# Original sheet is "Sheet1" (rId1/sheet1.xml). We need to insert a NEW sheet
# named "시트러스" BEFORE it, which becomes the new first tab and inherits all
# of Sheet1's current formatting/content (so copy it), then gets additional
# rows of data.  The old "Sheet1" is pushed to the second tab position and
# keeps its original 3-row content.
$orig = $wb.Worksheets.Item(1)
$orig.Copy($orig)

$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "시트러스"
$oldSheet = $wb.Worksheets.Item(2)

# --- Extend the new "시트러스" sheet with the additional rows/values ---
$newSheet.Range("B3").Value = "마바사아"
$newSheet.Range("B4").Value = "자차카타"
$newSheet.Range("B5").Value = "파하"
$newSheet.Range("C4").Value = "B12345124"
$newSheet.Range("C5").Value = "B23345125"
$newSheet.Range("D3").Value = 1600
$newSheet.Range("A4").Value = 1003
$newSheet.Range("D4").Value = 1800
$newSheet.Range("A5").Value = 1004
$newSheet.Range("D5").Value = 12000

# Print setup picked up by the new sheet.
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# Selection state: new sheet keeps focus (tabSelected) with a new active cell;
# old sheet's last selection is the original used range, and it is not the
# active tab.
[void]$oldSheet.Range("A1:D3").Select()
[void]$newSheet.Activate()
[void]$newSheet.Range("F24").Select()
